$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.66643366666667
$ws.Range("H2").Value = 37.999301
$ws.Range("I2").Value = 0.1759291503241684
$ws.Range("J2").Value = 0.1759291503241684
$ws.Range("M2").Value = 2.507757
$ws.Range("N2").Value = 7.523270999999999
$ws.Range("O2").Value = 0.07648041298707947
$ws.Range("P2").Value = 0.07648041298707947
$ws.Range("Q2").Value = 31.764337692619
$ws.Range("R2").Value = 285.879039233571
$ws.Range("S2").Value = 0.01345513407325839
$ws.Range("T2").Value = 0.01345513407325839
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.66643366666667
$ws.Range("H3").Value = 37.999301
$ws.Range("I3").Value = 0.1759291503241684
$ws.Range("J3").Value = 0.1759291503241684
$ws.Range("O3").Value = 0.6219651214303167
$ws.Range("P3").Value = 0.6219651214303167
$ws.Range("Q3").Value = 258.3185599884375
$ws.Range("R3").Value = 2324.867039895937
$ws.Range("S3").Value = 0.1094217953445038
$ws.Range("T3").Value = 0.1094217953445038
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.66643366666667
$ws.Range("H4").Value = 37.999301
$ws.Range("I4").Value = 0.1759291503241684
$ws.Range("J4").Value = 0.1759291503241684
$ws.Range("M4").Value = 9.887829999999999
$ws.Range("N4").Value = 29.66349
$ws.Range("O4").Value = 0.3015544655826039
$ws.Range("P4").Value = 0.301554465582604
$ws.Range("Q4").Value = 125.2435428022767
$ws.Range("R4").Value = 1127.19188522049
$ws.Range("S4").Value = 0.0530522209064062
$ws.Range("T4").Value = 0.0530522209064062
$ws.Range("I5").Value = 0.5164516272884614
$ws.Range("J5").Value = 0.5164516272884614
$ws.Range("M5").Value = 2.507757
$ws.Range("N5").Value = 7.523270999999999
$ws.Range("O5").Value = 0.07648041298707947
$ws.Range("P5").Value = 0.07648041298707947
$ws.Range("Q5").Value = 93.246308874145
$ws.Range("R5").Value = 839.2167798673049
$ws.Range("S5").Value = 0.03949843374287077
$ws.Range("T5").Value = 0.03949843374287077
$ws.Range("I6").Value = 0.5164516272884614
$ws.Range("J6").Value = 0.5164516272884614
$ws.Range("O6").Value = 0.6219651214303167
$ws.Range("P6").Value = 0.6219651214303167
$ws.Range("S6").Value = 0.3212148990793525
$ws.Range("T6").Value = 0.3212148990793525
$ws.Range("I7").Value = 0.5164516272884614
$ws.Range("J7").Value = 0.5164516272884614
$ws.Range("M7").Value = 9.887829999999999
$ws.Range("N7").Value = 29.66349
$ws.Range("O7").Value = 0.3015544655826039
$ws.Range("P7").Value = 0.301554465582604
$ws.Range("Q7").Value = 367.6606825442166
$ws.Range("R7").Value = 3308.94614289795
$ws.Range("S7").Value = 0.1557382944662381
$ws.Range("T7").Value = 0.1557382944662381
$ws.Range("G8").Value = 22.14777066666666
$ws.Range("H8").Value = 66.44331199999999
$ws.Range("I8").Value = 0.3076192223873702
$ws.Range("J8").Value = 0.3076192223873702
$ws.Range("M8").Value = 2.507757
$ws.Range("N8").Value = 7.523270999999999
$ws.Range("O8").Value = 0.07648041298707947
$ws.Range("P8").Value = 0.07648041298707947
$ws.Range("Q8").Value = 55.54122692372798
$ws.Range("R8").Value = 499.8710423135519
$ws.Range("S8").Value = 0.02352684517095032
$ws.Range("T8").Value = 0.02352684517095032
$ws.Range("G9").Value = 22.14777066666666
$ws.Range("H9").Value = 66.44331199999999
$ws.Range("I9").Value = 0.3076192223873702
$ws.Range("J9").Value = 0.3076192223873702
$ws.Range("O9").Value = 0.6219651214303167
$ws.Range("P9").Value = 0.6219651214303167
$ws.Range("Q9").Value = 451.6804316137937
$ws.Range("R9").Value = 4065.123884524144
$ws.Range("S9").Value = 0.1913284270064603
$ws.Range("T9").Value = 0.1913284270064603
$ws.Range("G10").Value = 22.14777066666666
$ws.Range("H10").Value = 66.44331199999999
$ws.Range("I10").Value = 0.3076192223873702
$ws.Range("J10").Value = 0.3076192223873702
$ws.Range("M10").Value = 9.887829999999999
$ws.Range("N10").Value = 29.66349
$ws.Range("O10").Value = 0.3015544655826039
$ws.Range("P10").Value = 0.301554465582604
$ws.Range("Q10").Value = 218.9933912309866
$ws.Range("R10").Value = 1970.94052107888
$ws.Range("S10").Value = 0.0927639502099596
$ws.Range("T10").Value = 0.09276395020995963

Write-Output "Updated cells"